# Apply "progress as of date 04 Nov 2025" update to the Training Dashboard sheet.
# For every data row (3-25): decrement the "PERIOD TO EXPIRE" (column H) value by 1
# and bump the "LAST UPDATE" (column I) date text from 03-Nov-2025 to 04-Nov-2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 25; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H - PERIOD TO EXPIRE
    $iCell = $ws.Cells.Item($row, 9)   # column I - LAST UPDATE

    $hCell.Value = $hCell.Value2 - 1

    # Prefix with an apostrophe so the date-looking text stays plain text
    # instead of being re-interpreted as a date serial number.
    $iCell.Value = "'04-Nov-2025"
}
